# Generate Report for Handoff
#
# This script updates the localization-status report to reflect that the
# files with Status "Ready for handoff" now have their handoff priority
# set to "ht" (matching the handoff type), and refreshes the associated
# handoff/report timestamps.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Rows (on every sheet) whose source file is in "Ready for handoff" status.
$rows = @(7, 8, 10, 11, 13, 14)

foreach ($r in $rows) {
    # Overview sheet: column G = "Latest HO Xliff Generate Date"
    $overview.Range("G$r").Value = "2016-09-06 14:37:48"

    # Column E = "Priority" -> should now match the handoff type "ht"
    $zhcn.Range("E$r").Value = "ht"
    $dede.Range("E$r").Value = "ht"

    # Column H = "Latest Handoff Datetime"
    $zhcn.Range("H$r").Value = "2016-09-06 14:37:41"
    $dede.Range("H$r").Value = "2016-09-06 14:37:48"
}
